$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 43: GMHO:0000260 "mental health intervention scenario plan" ---
# Inserting here shifts rows 43-71 down to 44-72, and copies the fill/style (s="3")
# from the row above (row 42), which matches the desired style for this new row.
$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value2 = "GMHO:0000260"
$ws.Range("B43").Value2 = "mental health intervention scenario plan"
$ws.Range("C43").Value2 = "A plan that is realized in a mental health intervention scenario process."
$ws.Range("D43").Value2 = "plan"
$ws.Range("Q43").Value2 = "LSR 1"
$ws.Range("T43").Value2 = "Published"

# --- Insert new row 57: OBI:0000260 "plan" ---
# After the first insert, "placebo intervention" is now row 56 and "plan specification"
# is now row 57. Insert a fresh row there to hold the new "plan" entry.
$ws.Rows.Item(57).Insert()

$ws.Range("A57").Value2 = "OBI:0000260"
$ws.Range("B57").Value2 = "plan"
$ws.Range("C57").Value2 = "A plan is a realizable entity that is the inheres in a bearer who is committed to realizing it as a planned process."
$ws.Range("D57").Value2 = "realizable"
$ws.Range("Q57").Value2 = "LSR 1"
$ws.Range("T57").Value2 = "External"

# Row 57 inherited the fill style of the row above it (style 2) via the Insert
# operation, but the source row has no fill at all - reset it back to Normal/no style.
$ws.Range("A57:W57").Style = "Normal"
